$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 716
$ws1.Range("F7").Value = 2730
$ws1.Range("F8").Value = 1661
$ws1.Range("F9").Value = 1735
$ws1.Range("F11").Value = 279
$ws1.Range("F12").Value = 708
$ws1.Range("F13").Value = 868
$ws1.Range("F14").Value = 145
$ws1.Range("F15").Value = 359
$ws1.Range("F18").Value = 45
$ws1.Range("F20").Value = 6228
$ws1.Range("F22").Value = 1316
$ws1.Range("F24").Value = 176
$ws1.Range("F26").Value = 292
$ws1.Range("F27").Value = 248
$ws1.Range("F28").Value = 55
$ws1.Range("F29").Value = 1085
$ws1.Range("F30").Value = 882
$ws1.Range("F34").Value = 450
$ws1.Range("F35").Value = 1299
$ws1.Range("F36").Value = 154
$ws1.Range("F38").Value = 205
$ws1.Range("F39").Value = 2
$ws1.Range("F40").Value = 137
$ws1.Range("F41").Value = 167
$ws1.Range("F42").Value = 138

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 14
$ws2.Range("F6").Value = 1

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 716
$ws4.Range("F6").Value = 14
$ws4.Range("F10").Value = 2730
$ws4.Range("F11").Value = 1662
$ws4.Range("F12").Value = 1735
$ws4.Range("F14").Value = 279
$ws4.Range("F15").Value = 708
$ws4.Range("F17").Value = 868
$ws4.Range("F18").Value = 145
$ws4.Range("F19").Value = 359
$ws4.Range("F21").Value = 45
$ws4.Range("F23").Value = 6228
$ws4.Range("F25").Value = 1316
$ws4.Range("F26").Value = 1
$ws4.Range("F30").Value = 292
$ws4.Range("F31").Value = 248
$ws4.Range("F32").Value = 55
$ws4.Range("F33").Value = 1085
$ws4.Range("F34").Value = 882
$ws4.Range("F38").Value = 450
$ws4.Range("F39").Value = 1299
$ws4.Range("F40").Value = 154
$ws4.Range("F42").Value = 205
$ws4.Range("F43").Value = 2
$ws4.Range("F44").Value = 137
$ws4.Range("F45").Value = 167
$ws4.Range("F49").Value = 138
